$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 101 (pushes the existing row 101 and everything
# below it down by one row).
$ws.Rows("101").Insert()

# Populate the newly inserted row 101 with the new record's data.
$ws.Range("A101").Value = 4
$ws.Range("B101").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C101").Value = "Los Lagos"
$ws.Range("D101").Value = 44781
$ws.Range("E101").Value = 10
$ws.Range("F101").Value = 100112009
$ws.Range("G101").Value = "Acelga"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 160
$ws.Range("K101").Value = 1200
$ws.Range("L101").Value = 1500
$ws.Range("M101").Value = 1350
$ws.Range("N101").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O101").Value = "Región de Los Lagos"
$ws.Range("P101").Value = 900
$ws.Range("Q101").Value = 1.5
$ws.Range("R101").Value = "Hortaliza"
